$d = $word.ActiveDocument
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_Hlk120140857"/><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Отчет по лабораторной работе 1</w:t></w:r></w:p><w:p><w:r><w:t>Цели</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r></w:p><w:p><w:r><w:t>Изучить цветовые модели: RGB, CMYK, HSV, HLS, XYZ, LAB, переход от одной модели к другой, исследовать цветовой график МКО. Создать приложение, позволяющее пользователю выбирать, а затем интерактивно менять цвет, показывая при этом его составляющие в нескольких моделях одновременно.</w:t></w:r></w:p><w:p><w:r><w:t>Ход работы</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r></w:p><w:p><w:r><w:t>-Предварительно сделан макет интерфейса и</w:t></w:r><w:r><w:t xml:space="preserve"> функций.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t>-И</w:t></w:r><w:r><w:t>зучены основные цветовые модели</w:t></w:r><w:r><w:t xml:space="preserve"> RGB, CMYK, HSV, HLS, XYZ, LAB</w:t></w:r><w:r><w:t xml:space="preserve"> и формулы всевозможных преобразований.</w:t></w:r></w:p><w:p><w:r><w:t>-</w:t></w:r><w:r><w:t xml:space="preserve">Изучены способы работы с цветами в </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>QT</w:t></w:r><w:r><w:t xml:space="preserve">, а также изучены возможности библиотеки </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>QPinter</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>-</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Сделан </w:t></w:r><w:r><w:t xml:space="preserve">интерфейс, отображающий все функции приложения, автоматический перевод цветов, отображение цветов на моделях, выбор цвета из палитры и с помощью конкретных значений. </w:t></w:r></w:p><w:p><w:r><w:t>-</w:t></w:r><w:r><w:t xml:space="preserve">Сделаны исключения в случае неправильных данных (часть ошибок обрабатывается самим </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>QT</w:t></w:r><w:r><w:t>).</w:t></w:r><w:r><w:t xml:space="preserve"> Перевод соответствует формулам и погрешность перевода незначительная. </w:t></w:r></w:p><w:p><w:r><w:t>-</w:t></w:r><w:r><w:t xml:space="preserve">Сделан </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>exe</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">файл, а исходные файлы с документацией загружены на </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>git</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t>Итог</w:t></w:r><w:r><w:t>:</w:t></w:r></w:p><w:p><w:r><w:t>Изуч</w:t></w:r><w:r><w:t>ены</w:t></w:r><w:r><w:t xml:space="preserve"> цветовые модели: RGB, CMYK, HSV, HLS, XYZ, LAB</w:t></w:r><w:r><w:t xml:space="preserve"> и способы </w:t></w:r><w:r><w:t xml:space="preserve">переход от одной модели к другой. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Изучены способы и особенности работы с цветами в </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Qt</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">Изучены классы </w:t></w:r><w:r><w:t>QColorDialog</w:t></w:r><w:r><w:t xml:space="preserve"> и </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>QColor</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>для работы с цветами.</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/></w:p><w:p><w:r><w:t>Созда</w:t></w:r><w:r><w:t>но</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>приложение</w:t></w:r><w:r><w:t xml:space="preserve"> выполняющее все поставленные цели.</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Content.InsertXML($xml)
